$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 57517
$ws.Range("B2").Value = "João Gabriel da Cunha"
$ws.Range("C2").Value = "Jurídico"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 45093
$ws.Range("G2").Value = 8559.34

# Row 3
$ws.Range("A3").Value = 73768
$ws.Range("B3").Value = "Pietro Nunes"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45079
$ws.Range("G3").Value = 3221.59

# Row 4
$ws.Range("A4").Value = 79404
$ws.Range("B4").Value = "Emanuel da Luz"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45079
$ws.Range("G4").Value = 10637.93

# Row 5
$ws.Range("A5").Value = 12781
$ws.Range("B5").Value = "Caroline Santos"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45086
$ws.Range("G5").Value = 10271.26

# Row 6
$ws.Range("A6").Value = 91636
$ws.Range("B6").Value = "Gustavo da Luz"
$ws.Range("C6").Value = "Atendimento ao Cliente"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45086
$ws.Range("G6").Value = 4808.21

# Row 7
$ws.Range("A7").Value = 82423
$ws.Range("B7").Value = "Felipe Novaes"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45086
$ws.Range("G7").Value = 5433.92

# Row 8
$ws.Range("A8").Value = 89273
$ws.Range("B8").Value = "Marcelo da Costa"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45078
$ws.Range("G8").Value = 5637.72

# Row 9
$ws.Range("A9").Value = 29557
$ws.Range("B9").Value = "Maria Ferreira"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45080
$ws.Range("G9").Value = 9908.370000000001

# Row 10
$ws.Range("A10").Value = 43722
$ws.Range("B10").Value = "Sr. Danilo da Rosa"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45100
$ws.Range("G10").Value = 7945.46

# Row 11
$ws.Range("A11").Value = 51391
$ws.Range("B11").Value = "Maria Fernanda Novaes"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 3245.58
